$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50")
foreach ($cellName in $textCells) {
    $ws.Range($cellName).NumberFormat = "@"
}

$ws.Range("D2").Value = "245.77"
$ws.Range("D3").Value = "24.21"
$ws.Range("D4").Value = "5.284"
$ws.Range("D5").Value = "0.05770"
$ws.Range("D6").Value = "6.490"
$ws.Range("D7").Value = "3.147"
$ws.Range("D8").Value = "0.8169"
$ws.Range("D9").Value = "0.8578"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.009752"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1363"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.06951"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03138"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.02897"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09392"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.753"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001525"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04663"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "0.006105"
$ws.Range("D20").Value = "0.001238"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("D22").Value = "0.00006104"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("D23").Value = "3.500"
$ws.Range("D24").Value = "2.137"
$ws.Range("D25").Value = "0.3195"
$ws.Range("D26").Value = "0.1357"
$ws.Range("D27").Value = "0.1329"
$ws.Range("D28").Value = "0.0002333"
$ws.Range("D40").Value = "0.03677"
$ws.Range("D41").Value = "0.006264"
$ws.Range("D42").Value = "0.1054"
$ws.Range("D43").Value = "0.002802"
$ws.Range("D44").Value = "0.008496"
$ws.Range("D45").Value = "0.00005276"
$ws.Range("D47").Value = "0.3703"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "0.002321"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").Value = "0.0002001"
